# Add a new "Failing login" test case as row 5 (sheet row 7) of the
# "To Do App (Navigation and Login)" test-scenario table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (test case #4) already carries the exact formatting we need for the
# new row (number style, wrap-text body style, vertical-centered outcome
# style) - copy it down to row 7 so the new cells inherit identical styles
# instead of the generic placeholder style that was there before.
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new test case contents (entered in the same left-to-right,
# "name -> description -> expected outcome -> values" order a person typing
# across the row would use).
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "test_<LoginIncorrectUsernameandPassword>"
$ws.Range("C7").Value = "This is to test whether user can login with incorrect username or password"
$ws.Range("E7").Value = "Error is displayed as such: ""You have entered the wrong username/password. Please try again!"""
$ws.Range("D7").Value = "username: ""username""" + "`n" + "password: ""pssword""" + "`n" + "username: ""wronguser""" + "`n" + "password: ""p@ssw0rd"""

# Entering multi-line text auto-grows the row; the sheet uses a fixed
# (customHeight) 15.75pt row height throughout the table, so restore it.
$ws.Range("A7").EntireRow.RowHeight = 15.75

# Leave the new "Expected Outcome" cell of the freshly entered row selected,
# matching where the author's cursor ended up.
$ws.Range("E7").Select()
